# Auto-generated edit script for cs-en-us-121pct.xlsx weekly refresh
# 'New crime data collected' - updates report header (volume/week dates)
# and refreshes the crime-complaint statistics table (rows 15-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number and report week dates ---
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Crime complaint statistics table ---
$ws.Range("M14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Formula = "=""0"""
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("M14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Formula = "=""***.*"""
$ws.Range("E15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4163)
$ws.Range("C16").Value = 2
$ws.Range("M14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Formula = "=""0"""
$ws.Range("D16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("M14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Formula = "=""***.*"""
$ws.Range("E16").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4163)
$ws.Range("F16").Value = 8
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 77
$ws.Range("K16").Value = 26.229508196721
$ws.Range("L16").Value = 1.315789473684
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -62.5
$ws.Range("F17").Value = 14
$ws.Range("H17").Value = -30
$ws.Range("I17").Value = 203
$ws.Range("J17").Value = 166
$ws.Range("K17").Value = 22.289156626506
$ws.Range("L17").Value = 11.538461538461
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 85
$ws.Range("J18").Value = 78
$ws.Range("K18").Value = 8.974358974358
$ws.Range("L18").Value = 8.974358974358
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -37.5
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -6.666666666666
$ws.Range("I19").Value = 380
$ws.Range("J19").Value = 304
$ws.Range("K19").Value = 25
$ws.Range("L19").Value = 33.802816901408
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 8
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = 112
$ws.Range("L20").Value = 103.846153846154
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -26.315789473684
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = 5.970149253731
$ws.Range("I21").Value = 864
$ws.Range("J21").Value = 680
$ws.Range("K21").Value = 27.058823529411
$ws.Range("L21").Value = 24.137931034482
$ws.Range("M14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Formula = "=""0"""
$ws.Range("C23").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4163)
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -12.5
$ws.Range("F24").Value = 110
$ws.Range("G24").Value = 100
$ws.Range("H24").Value = 10
$ws.Range("I24").Value = 1331
$ws.Range("J24").Value = 901
$ws.Range("K24").Value = 47.724750277469
$ws.Range("L24").Value = 62.713936430317
$ws.Range("C25").Value = 17
$ws.Range("E25").Value = 54.545454545454
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -12.195121951219
$ws.Range("I25").Value = 471
$ws.Range("J25").Value = 396
$ws.Range("K25").Value = 18.939393939393
$ws.Range("L25").Value = 41.017964071856
$ws.Range("M14").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Formula = "=""0"""
$ws.Range("C26").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4163)
$ws.Range("M14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Formula = "=""0"""
$ws.Range("D26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("M14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("E26").Formula = "=""***.*"""
$ws.Range("E26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 100
$ws.Range("L26").Value = -12.5
$ws.Range("F14").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 100
$ws.Range("F14").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 40
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -13.043478260869
$ws.Range("L27").Value = 73.913043478260
$ws.Range("F14").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = 2
$ws.Range("M14").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Formula = "=""0"""
$ws.Range("G30").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4163)
$ws.Range("M14").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Formula = "=""***.*"""
$ws.Range("H30").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("I30").Value = 10
$ws.Range("K30").Value = 233.333333333333
$ws.Range("L30").Value = 900

$excel.CutCopyMode = $false
